$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 739.125
$ws.Range("I11").Value = 739.125
$ws.Range("K11").Value = 739.125
$ws.Range("M11").Value = -599.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1926.6666
$ws.Range("J17").Value = 1926.6666
$ws.Range("L17").Value = 5779.9998
$ws.Range("N17").Value = -6115.9998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1199.7142
$ws.Range("I101").Value = 1279.8
$ws.Range("K101").Value = 3839.4
$ws.Range("M101").Value = -2217.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4144.1665
$ws.Range("J116").Value = 3434.5
$ws.Range("L116").Value = 3434.5
$ws.Range("N116").Value = -10318.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4962.067
$ws.Range("J138").Value = 7803.4546
$ws.Range("L138").Value = 23410.3638
$ws.Range("N138").Value = -33690.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2568.4443
$ws.Range("I2").Value = 1902.7142
$ws.Range("K2").Value = 1902.7142
$ws.Range("M2").Value = -1789.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6533.316
$ws.Range("I32").Value = 3539.7812
$ws.Range("K32").Value = 3539.7812
$ws.Range("M32").Value = -3252.7812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2568.4443
$ws.Range("I116").Value = 1902.7142
$ws.Range("K116").Value = 1902.7142
$ws.Range("M116").Value = 391.2858000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1041.1111
$ws.Range("I132").Value = 1012.9583
$ws.Range("J132").Value = 1266.3334
$ws.Range("K132").Value = 3038.8749
$ws.Range("L132").Value = 3799.0002
$ws.Range("M132").Value = -508.8748999999998
$ws.Range("N132").Value = -8859.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2568.4443
$ws.Range("I3").Value = 1902.7142
$ws.Range("K3").Value = 1902.7142
$ws.Range("M3").Value = -1788.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5427.909
$ws.Range("I20").Value = 4687.25
$ws.Range("K20").Value = 4687.25
$ws.Range("M20").Value = -4440.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4102
$ws.Range("I105").Value = 4195.4
$ws.Range("J105").Value = 3635
$ws.Range("K105").Value = 4195.4
$ws.Range("L105").Value = 3635
$ws.Range("M105").Value = -2448.4
$ws.Range("N105").Value = -7129

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2038
$ws.Range("I134").Value = 1915.4375
$ws.Range("K134").Value = 5746.3125
$ws.Range("M134").Value = -3211.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1000
$ws.Range("J14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("N14").Value = -1340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5511.278
$ws.Range("I31").Value = 4947.75
$ws.Range("K31").Value = 4947.75
$ws.Range("M31").Value = -4652.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5511.278
$ws.Range("I34").Value = 4947.75
$ws.Range("K34").Value = 4947.75
$ws.Range("M34").Value = -4745.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11453.704
$ws.Range("I99").Value = 8272.125
$ws.Range("K99").Value = 8272.125
$ws.Range("M99").Value = -6774.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 11453.704
$ws.Range("I126").Value = 8272.125
$ws.Range("K126").Value = 24816.375
$ws.Range("M126").Value = -22346.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 634.4
$ws.Range("J17").Value = 1039
$ws.Range("L17").Value = 3117
$ws.Range("N17").Value = -3455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 171.36363
$ws.Range("I23").Value = 187.14285
$ws.Range("K23").Value = 561.4285500000001
$ws.Range("M23").Value = -326.4285500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 2168.375
$ws.Range("I24").Value = 316
$ws.Range("J24").Value = 3279.8
$ws.Range("K24").Value = 948
$ws.Range("L24").Value = 9839.400000000001
$ws.Range("M24").Value = -718
$ws.Range("N24").Value = -10299.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 341.33334
$ws.Range("I35").Value = 345.33334
$ws.Range("J35").Value = 333.33334
$ws.Range("K35").Value = 1036.00002
$ws.Range("L35").Value = 1000.00002
$ws.Range("M35").Value = -748.0000199999999
$ws.Range("N35").Value = -1576.00002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 106.8421
$ws.Range("J40").Value = 114.77778
$ws.Range("L40").Value = 459.11112
$ws.Range("N40").Value = -597.11112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1268.4
$ws.Range("J131").Value = 1460.0667
$ws.Range("L131").Value = 4380.2001
$ws.Range("N131").Value = -14460.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4529
$ws.Range("I80").Value = 4529
$ws.Range("K80").Value = 4529
$ws.Range("M80").Value = -3531

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4529
$ws.Range("I83").Value = 4529
$ws.Range("K83").Value = 22645
$ws.Range("M83").Value = -17653

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2300.9285
$ws.Range("I102").Value = 1337.1666
$ws.Range("K102").Value = 1337.1666
$ws.Range("M102").Value = 284.8334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2723.4119
$ws.Range("I132").Value = 2230.6
$ws.Range("K132").Value = 6691.799999999999
$ws.Range("M132").Value = -4161.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3465.5
$ws.Range("I7").Value = 3465.5
$ws.Range("K7").Value = 3465.5
$ws.Range("M7").Value = -3353.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2898.8
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1206.6
$ws.Range("I61").Value = 1133.25
$ws.Range("K61").Value = 1133.25
$ws.Range("M61").Value = -931.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1206.6
$ws.Range("I113").Value = 1133.25
$ws.Range("K113").Value = 1133.25
$ws.Range("M113").Value = 1036.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5699.6
$ws.Range("I122").Value = 5500
$ws.Range("K122").Value = 16500
$ws.Range("M122").Value = -14050

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3465.5
$ws.Range("I126").Value = 3465.5
$ws.Range("K126").Value = 10396.5
$ws.Range("M126").Value = -7926.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5254.4
$ws.Range("I132").Value = 4029.8
$ws.Range("K132").Value = 12089.4
$ws.Range("M132").Value = -9559.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7430.6523
$ws.Range("I81").Value = 1750.4
$ws.Range("J81").Value = 11800.077
$ws.Range("K81").Value = 3500.8
$ws.Range("L81").Value = 23600.154
$ws.Range("M81").Value = -2439.8
$ws.Range("N81").Value = -25722.154

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7430.6523
$ws.Range("I84").Value = 1750.4
$ws.Range("J84").Value = 11800.077
$ws.Range("K84").Value = 17504
$ws.Range("L84").Value = 118000.77
$ws.Range("M84").Value = -12200
$ws.Range("N84").Value = -128608.77

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1049.8
$ws.Range("I96").Value = 1055.4445
$ws.Range("J96").Value = 999
$ws.Range("K96").Value = 1055.4445
$ws.Range("L96").Value = 999
$ws.Range("M96").Value = 317.5554999999999
$ws.Range("N96").Value = -3745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1776
$ws.Range("I122").Value = 1035.4
$ws.Range("J122").Value = 3997.8
$ws.Range("K122").Value = 3106.2
$ws.Range("L122").Value = 11993.4
$ws.Range("M122").Value = -656.2000000000003
$ws.Range("N122").Value = -16893.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2710.261
$ws.Range("J126").Value = 5638.5713
$ws.Range("L126").Value = 16915.7139
$ws.Range("N126").Value = -21855.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1051.4138
$ws.Range("I136").Value = 1071.5
$ws.Range("K136").Value = 3214.5
$ws.Range("M136").Value = -664.5
